$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Aton picked 3 more Tree questions - append as rows 44-46, re-using the
# formatting already used for similar "Easy"/"Medium" rows further up
# (row 41 for an "Easy" row, row 43 for a "Medium" row), but skipping
# column E entirely since these new rows have no status/comment cell.

# --- Row 44: Binary Tree Level Order Traversal II (Easy) ---
$ws.Range("A41:D41").Copy()
$ws.Range("A44:D44").PasteSpecial(-4122)
$ws.Range("F41:G41").Copy()
$ws.Range("F44:G44").PasteSpecial(-4122)

$ws.Range("A44").Value = 107
$ws.Range("B44").Value = "Binary Tree Level Order Traversal II"
$ws.Range("C44").Value = "Tree"
$ws.Range("D44").Value = "Aton"
$ws.Range("F44").Value = "Easy"
$ws.Range("G44").Value = "Python"

# --- Row 45: Maximum Binary Tree (Medium) ---
$ws.Range("A43:D43").Copy()
$ws.Range("A45:D45").PasteSpecial(-4122)
$ws.Range("F43:G43").Copy()
$ws.Range("F45:G45").PasteSpecial(-4122)

$ws.Range("A45").Value = 654
$ws.Range("B45").Value = "Maximum Binary Tree"
$ws.Range("C45").Value = "Tree"
$ws.Range("D45").Value = "Aton"
$ws.Range("F45").Value = "Medium"
$ws.Range("G45").Value = "Python"

# --- Row 46: Two Sum IV - Input is a BST (Easy) ---
$ws.Range("A41:D41").Copy()
$ws.Range("A46:D46").PasteSpecial(-4122)
$ws.Range("F41:G41").Copy()
$ws.Range("F46:G46").PasteSpecial(-4122)

$ws.Range("A46").Value = 653
$ws.Range("B46").Value = "Two Sum IV - Input is a BST"
$ws.Range("C46").Value = "Tree"
$ws.Range("D46").Value = "Aton"
$ws.Range("F46").Value = "Easy"
$ws.Range("G46").Value = "Python"

$excel.CutCopyMode = $false

# Match the updated selection left behind in the saved file
$ws.Range("H45").Select()
